# 2025-08 "ooutput update" refresh:
#  - repoint the published URLs from the old GitHub/shorthand location to
#    the new 2rdoc.pt IG host (StructureDefinition + ValueSet canonicals)
#  - bump the generation Date metadata value
#  - column widths on the "Elements" sheet shift because several best-fit
#    columns are driven by the (now shorter) canonical URL text

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet: URL + Date --------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/sleep-quality"
$wsMeta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# ---- Elements sheet: Binding Value Set URL ----------------------------
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("Z6").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ValueSet/sleep-quality-extended-vs"

# Extension.url's "Fixed Value" cell (R5) duplicates the same canonical
# StructureDefinition URL shown on the Metadata sheet (same shared string
# in the source workbook), so it moves together with B2 above.
$wsElem.Range("R5").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/sleep-quality"

# ---- Elements sheet: best-fit column widths ---------------------------
# (target widths taken from the canonical OOXML diff; ColumnWidth is
# expressed in character units, quantized to the nearest 1/6 by the host,
# so each value below is chosen to land on the closest achievable width.)
$wsElem.Columns.Item(1).ColumnWidth = 15.666666666666666
$wsElem.Columns.Item(2).ColumnWidth = 15.666666666666666
$wsElem.Columns.Item(3).ColumnWidth = 9.0
$wsElem.Columns.Item(3).Hidden = $true
$wsElem.Columns.Item(4).ColumnWidth = 6.166666666666667
$wsElem.Columns.Item(4).Hidden = $true
$wsElem.Columns.Item(5).ColumnWidth = 4.5
$wsElem.Columns.Item(6).ColumnWidth = 3.1666666666666665
$wsElem.Columns.Item(7).ColumnWidth = 3.5
$wsElem.Columns.Item(8).ColumnWidth = 11.833333333333334
$wsElem.Columns.Item(9).ColumnWidth = 9.666666666666666
$wsElem.Columns.Item(11).ColumnWidth = 13.5
$wsElem.Columns.Item(15).ColumnWidth = 11.5
$wsElem.Columns.Item(20).ColumnWidth = 7.0
$wsElem.Columns.Item(21).ColumnWidth = 12.833333333333334
$wsElem.Columns.Item(22).ColumnWidth = 13.166666666666666
$wsElem.Columns.Item(23).ColumnWidth = 14.166666666666666
$wsElem.Columns.Item(24).ColumnWidth = 13.833333333333334
$wsElem.Columns.Item(25).ColumnWidth = 16.166666666666668
$wsElem.Columns.Item(26).ColumnWidth = 61.666666666666664
$wsElem.Columns.Item(27).ColumnWidth = 4.166666666666667
$wsElem.Columns.Item(28).ColumnWidth = 17.166666666666668
$wsElem.Columns.Item(29).ColumnWidth = 33.666666666666664
$wsElem.Columns.Item(30).ColumnWidth = 12.666666666666666
$wsElem.Columns.Item(31).ColumnWidth = 10.5
$wsElem.Columns.Item(31).Hidden = $true
$wsElem.Columns.Item(32).ColumnWidth = 14.166666666666666
$wsElem.Columns.Item(32).Hidden = $true
$wsElem.Columns.Item(33).ColumnWidth = 7.333333333333333
$wsElem.Columns.Item(33).Hidden = $true
$wsElem.Columns.Item(34).ColumnWidth = 7.666666666666667
$wsElem.Columns.Item(37).ColumnWidth = 18.666666666666668
